$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab from "Through 2022-03-22" to "Through 2022-03-23"
$ws.Name = "Through 2022-03-23"

# Update the "March" row label text
$ws.Range("A4").Value = "March (through 03-23)"

# Update the March row (row 4) values
$ws.Range("B4").Value = 23
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 47
$ws.Range("F4").Value = 22
$ws.Range("G4").Value = 44
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 97

# Update the Total row (row 5) values
$ws.Range("B5").Value = 60
$ws.Range("D5").Value = 171
$ws.Range("E5").Value = 184
$ws.Range("F5").Value = 101
$ws.Range("G5").Value = 185
$ws.Range("H5").Value = 402
$ws.Range("I5").Value = 397
